# Update Name of Algo
# Applies updated KNN-imputed values in column E for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 17.332
    13 = 16.498
    16 = 16.696
    18 = 16.626
    20 = 16.314
    26 = 16.424
    27 = 16.39
    29 = 16.85
    35 = 16.407
    36 = 16.69
    45 = 16.804
    55 = 16.557
    57 = 16.564
    69 = 17.4
    76 = 16.433
    78 = 16.507
    82 = 16.864
    83 = 16.935
    93 = 17.291
    97 = 16.86
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
